$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete "Markets and Economy" (row 19) and "Moving resources" (row 20)
# passive-skill rows entirely - they're being replaced by a single new
# "Market Place" building/skill described in row 18.
$ws.Rows("19:20").Delete()

# Row 18 ("Capital city") becomes the new "Market Place" entry, with an
# updated description and tweaked effect_type / unlocks_at_level values.
$ws.Range("B18").Value = "Market Place"
$ws.Range("C18").Value = "Allows players to build and upgrade the Market Place building for their kingdoms. Each kingdom that has a Market Place leveled to level 5, can then request resources from other kingdoms you own, on the same plane - who also have a Market Place at level 5."
$ws.Range("H18").Value = 4
$ws.Range("I18").Value = 15
